$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Data edits (E0 row / row 7: B7, D7) ---
$ws.Range("B7").Value = 2140
$ws.Range("D7").Value = 98

# --- Data edits (E2 row / row 9: D9) ---
$ws.Range("D9").Value = 97.8

# --- Data edits (E3 row / row 10: D10) ---
$ws.Range("D10").Value = 103

# --- Re-apply center alignment to the "Actual Distance Moved" column for
#     E1/E2/E3 rows (D8:D10) - matches the author re-touching that formatting ---
$ws.Range("D8:D10").HorizontalAlignment = -4108

# --- Update the active selection left behind by the author (D8 -> D11) ---
$ws.Range("D11").Select()
